# Insert a new daily price record for Pomelo (Feria Lagunitas de Puerto Montt)
# at row 57, pushing the existing rows 57-176 down to 58-177.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(57).Insert()

$ws.Cells.Item(57, 1).Value = 4
$ws.Cells.Item(57, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value = "Los Lagos"
$ws.Cells.Item(57, 4).Value = 44519
$ws.Cells.Item(57, 5).Value = 10
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100102
$ws.Cells.Item(57, 8).Value = "Cítricos"
$ws.Cells.Item(57, 9).Value = 100102006
$ws.Cells.Item(57, 10).Value = "Pomelo"
$ws.Cells.Item(57, 11).Value = "Start Ruby"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 200
$ws.Cells.Item(57, 14).Value = 11000
$ws.Cells.Item(57, 15).Value = 12000
$ws.Cells.Item(57, 16).Value = 11500
$ws.Cells.Item(57, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(57, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(57, 19).Value = 821
$ws.Cells.Item(57, 20).Value = 14
